$d = $word.ActiveDocument

# 1. Increase the header-row height of the two "detailed" result tables
#    (the ones with the chi-squared column) from 571 twips (28.55 pt)
#    to 637 twips (31.85 pt) to reflect the re-run with 1000 iterations.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $r = $t.Rows.Item(1)
    if ($r.Height -eq 28.55) {
        $r.Height = 31.85
    }
}

# 2. Fix the mis-encoded chi (χ) character in the "χ²" column headers of
#    those same tables.
$d.Content.Find.Execute("χ", $false, $false, $false, $false, $false, $true, 1, $false, "Ï‡", 2)
